# Updates data.xlsx per the "Cập nhật data.xlsx từ công cụ QR" commit:
#  - rename header "group" -> "branch"
#  - insert a new data row (QR-tool generated record) above the existing sample row
#  - drop the styled/bold header look (back to plain default style)
#  - normalize column widths, unhide the two previously-hidden trailing columns
#  - remove the autoFilter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header: "group" -> "branch" ---------------------------------------
$ws.Range("I1").Value = "branch"

# --- insert the new record as row 2, pushing the sample row to row 3 ---
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "a4do63e3aba"
$ws.Range("B2").Value = "TESTM9R7"
$ws.Range("C2").Value = "Hộ kinh doanh Test 85H7"
$ws.Range("D2").Value = "02 Hòa Bình, Ninh Kiều, Cần Thơ"
$ws.Range("E2").Value = "https://www.google.com/maps/search/?api=1&query=02%20H%C3%B2a%20B%C3%ACnh%2C%20Ninh%20Ki%E1%BB%81u%2C%20C%E1%BA%A7n%20Th%C6%A1"
$ws.Range("F2").Value = "2025-08-13T08:49:42.394Z"
$ws.Range("G2").Value = "Dòng test thêm nhanh"
# leading apostrophe forces text so the leading zero in these ID-like
# numbers survives the round-trip instead of becoming a number
$ws.Range("H2").Value = "'0908405063"
$ws.Range("I2").Value = "CN Cần Thơ II"
$ws.Range("J2").Value = "'0405341183179"
$ws.Range("K2").Value = "KH3284"
$ws.Range("L2").Value = "Demo User"
# M2 (pinSalt) / N2 (pinHash) stay blank, matching the source row's empty values

# the inserted row inherits the header's bold/shaded look - strip it back
# down to the workbook's plain default style
$ws.Rows.Item(2).ClearFormats()

# --- drop the bold/shaded header styling ---------------------------------
$ws.Rows.Item(1).ClearFormats()
$ws.Rows.Item(1).AutoFit()

# --- column widths: unhide M:N, then normalize widths ---------------------
$ws.Range("A1:N1").EntireColumn.Hidden = $false
$ws.Range("A1:B1").EntireColumn.ColumnWidth = 17.1
$ws.Range("F1:N1").EntireColumn.ColumnWidth = 17.1

# --- remove the autofilter ------------------------------------------------
$ws.AutoFilterMode = $false
